# Updated cryptos list on Fri Feb 17 20:33:48 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.692.95"
$ws.Range("E2").Value = "'  +0.26%  "
$ws.Range("D3").Value = "'1.702.35"
$ws.Range("E3").Value = "'  +0.47%  "
$ws.Range("D4").Value = "'0.9943"
$ws.Range("E4").Value = "'  -1.41%  "
$ws.Range("D5").Value = "'312.98"
$ws.Range("E5").Value = "'  -1.56%  "
$ws.Range("D6").Value = "'0.9927"
$ws.Range("E6").Value = "'  -1.23%  "
$ws.Range("D7").Value = "'0.3949"
$ws.Range("E7").Value = "'  -0.49%  "
$ws.Range("D8").Value = "'0.4065"
$ws.Range("E8").Value = "'  +1.61%  "
$ws.Range("D9").Value = "'1.517"
$ws.Range("E9").Value = "'  +6.61%  "
$ws.Range("D10").Value = "'0.9946"
$ws.Range("E10").Value = "'  -1.47%  "
$ws.Range("D11").Value = "'53.43"
$ws.Range("E11").Value = "'  +9.82%  "
$ws.Range("D12").Value = "'0.08768"
$ws.Range("E12").Value = "'  -0.33%  "
$ws.Range("D13").Value = "'7.302"
$ws.Range("E13").Value = "'  +10.07%  "
$ws.Range("D14").Value = "'23.28"
$ws.Range("E14").Value = "'  +0.07%  "
$ws.Range("D15").Value = "'0.00001322"
$ws.Range("E15").Value = "'  -0.03%  "
$ws.Range("D16").Value = "'7.468"
$ws.Range("E16").Value = "'  +2.84%  "
$ws.Range("D17").Value = "'1.699.09"
$ws.Range("E17").Value = "'  -0.55%  "
$ws.Range("D18").Value = "'100.63"
$ws.Range("E18").Value = "'  -1.76%  "
$ws.Range("D19").Value = "'0.07017"
$ws.Range("E19").Value = "'  +2.16%  "
$ws.Range("D20").Value = "'19.48"
$ws.Range("E20").Value = "'  -0.75%  "
$ws.Range("D21").Value = "'6.748"
$ws.Range("E21").Value = "'  -1.20%  "
$ws.Range("D22").Value = "'0.9928"
$ws.Range("E22").Value = "'  -1.26%  "
$ws.Range("E23").Value = "'  +1.61%  "
$ws.Range("D24").Value = "'24.695.11"
$ws.Range("E24").Value = "'  +0.22%  "
$ws.Range("D25").Value = "'2.962"
$ws.Range("E25").Value = "'  +2.73%  "
$ws.Range("D26").Value = "'2.304"
$ws.Range("E26").Value = "'  -0.68%  "
$ws.Range("D27").Value = "'22.39"
$ws.Range("E27").Value = "'  +0.48%  "
$ws.Range("D28").Value = "'158.48"
$ws.Range("E28").Value = "'  -1.34%  "
$ws.Range("D29").Value = "'5.122"
$ws.Range("E29").Value = "'  -3.29%  "
$ws.Range("D30").Value = "'133.22"
$ws.Range("E30").Value = "'  -0.08%  "
$ws.Range("D31").Value = "'7.426"
$ws.Range("E31").Value = "'  +25.74%  "
$ws.Range("B32").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "'1.883.97"
$ws.Range("E32").Value = "'  -0.69%  "
$ws.Range("B33").Value = "'ImmutableX"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.094"
$ws.Range("E33").Value = "'  -8.57%  "
$ws.Range("B34").Value = "'InternetComputer(DFINITY)"
$ws.Range("C34").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'7.411"
$ws.Range("E34").Value = "'  +21.23%  "
$ws.Range("B35").Value = "'Hedera"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.08688"
$ws.Range("E35").Value = "'  -4.30%  "
$ws.Range("D36").Value = "'11.14"
$ws.Range("E36").Value = "'  +1.08%  "
$ws.Range("D37").Value = "'0.2732"
$ws.Range("E37").Value = "'  +0.69%  "
$ws.Range("B38").Value = "'WEMIXTOKEN"
$ws.Range("C38").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.934"
$ws.Range("E38").Value = "'  +3.15%  "
$ws.Range("B39").Value = "'Aptos"
$ws.Range("C39").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'14.78"
$ws.Range("E39").Value = "'  -4.65%  "
$ws.Range("D40").Value = "'0.02765"
$ws.Range("E40").Value = "'  +8.94%  "
$ws.Range("D41").Value = "'0.08960"
$ws.Range("E41").Value = "'  +0.38%  "
$ws.Range("D42").Value = "'1.473"
$ws.Range("E42").Value = "'  +0.17%  "
$ws.Range("D43").Value = "'0.7657"
$ws.Range("E43").Value = "'  +0.68%  "
$ws.Range("D44").Value = "'0.7223"
$ws.Range("E44").Value = "'  +0.75%  "
$ws.Range("D45").Value = "'15.36"
$ws.Range("E45").Value = "'  +0.39%  "
$ws.Range("D46").Value = "'2.455"
$ws.Range("E46").Value = "'  -0.48%  "
$ws.Range("D47").Value = "'4.154"
$ws.Range("E47").Value = "'  -0.11%  "
$ws.Range("D48").Value = "'0.9925"
$ws.Range("E48").Value = "'  -1.25%  "
$ws.Range("D49").Value = "'141.84"
$ws.Range("E49").Value = "'  -0.78%  "
$ws.Range("D50").Value = "'1.310"
$ws.Range("E50").Value = "'  +12.80%  "
$ws.Range("D51").Value = "'0.08033"
$ws.Range("E51").Value = "'  +0.98%  "
